# Update "想去人数" (number of people interested) counts on the
# "展览" and "全部类型" sheets to reflect the latest scrape.
#   F2: 607  -> 612
#   F3: 3740 -> 3755
#   F5: 717  -> 718

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 612
    $ws.Range("F3").Value = 3755
    $ws.Range("F5").Value = 718
}
